# Updates cryptos list values (Price column D, Volume(1h) column E)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.167.29"
$ws.Range("E2").Value = "  +0.47%  "
$ws.Range("D3").Value = "1.684.36"
$ws.Range("E3").Value = "  +0.08%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.90"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  -0.11%  "
$ws.Range("E6").Value = "  +0.18%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "22.86"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +5.83%  "
$ws.Range("E9").Value = "  +2.52%  "
$ws.Range("E10").Value = "  +0.34%  "
$ws.Range("E11").Value = "  +0.00%  "
$ws.Range("D12").Value = "1.923.55"
$ws.Range("E12").Value = "  +0.18%  "
$ws.Range("D13").Value = "1.688.95"
$ws.Range("E13").Value = "  +1.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "4.19"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +2.09%  "
$ws.Range("E15").Value = "  +4.39%  "
$ws.Range("E16").Value = "  +0.80%  "
$ws.Range("D17").Value = "27.185.40"
$ws.Range("E17").Value = "  +0.52%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "235.59"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  -0.46%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.99"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -2.42%  "
$ws.Range("D20").Value = "0.0₃0742"
$ws.Range("E20").Value = "  +0.68%  "
$ws.Range("E21").Value = "  +0.04%  "
$ws.Range("E22").Value = "  +2.10%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "9.57"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +3.08%  "
$ws.Range("E24").Value = "  -1.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "146.91"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E26").Value = "  +1.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.42"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -2.25%  "
$ws.Range("E28").Value = "  +0.02%  "
$ws.Range("E29").Value = "  +0.16%  "
$ws.Range("E30").Value = "  +0.91%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.17"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -0.10%  "
$ws.Range("E32").Value = "  +0.83%  "
$ws.Range("D33").Value = "1.544.36"
$ws.Range("E33").Value = "  +1.40%  "
$ws.Range("E34").Value = "  +1.77%  "
$ws.Range("E35").Value = "  -2.06%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.604"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +2.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.945"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +2.58%  "
$ws.Range("E38").Value = "  -0.60%  "
$ws.Range("E39").Value = "  -1.27%  "
$ws.Range("E40").Value = "  +2.25%  "
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "69.08"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.40%  "
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("D45").Value = "1.830.23"
$ws.Range("E45").Value = "  +0.41%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.791"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +1.30%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "90.14"
$ws.Range("D47").ClearFormats()
$ws.Range("E48").Value = "  +6.19%  "
$ws.Range("E49").Value = "  +5.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.23"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +3.29%  "
